$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF ("Date") was stamped with the workbook's own file-name fragment
# ("6-4-2007-08") on every data row instead of the real game date. NBA.com's
# box-score export showed the date a day off (timezone rollover), so the
# correct date for this file is 2008-06-04. Fix rows 2-31 (the full data
# range) and format the range as Text first so Excel keeps the corrected
# value as a literal string instead of reinterpreting the date-shaped text
# as a serial date number.
$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"
$rng.Value = "2008-06-04"
